# Refresh crypto market data (price + 1h volume change) for cryptos.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: numeric-looking price values are written with a leading apostrophe
# so Excel keeps them as text, matching the original inline-string cells
# (e.g. "0.9998") instead of silently converting them to numbers.

$ws.Range('D2').Value = '30.023.88'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '1.909.33'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('D4').Value = '''0.9998'
$ws.Range('D5').Value = '''0.8291'
$ws.Range('E5').Value = '  +8.82%  '
$ws.Range('D6').Value = '''242.09'
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('D7').Value = '''0.9998'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '''0.3224'
$ws.Range('E8').Value = '  +5.32%  '
$ws.Range('D9').Value = '''26.68'
$ws.Range('E9').Value = '  +4.25%  '
$ws.Range('D10').Value = '''0.07011'
$ws.Range('E10').Value = '  +2.48%  '
$ws.Range('D11').Value = '''0.08019'
$ws.Range('E11').Value = '  +0.73%  '
$ws.Range('D12').Value = '''0.7499'
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('D13').Value = '1.912.82'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').Value = '''5.211'
$ws.Range('E14').Value = '  +0.75%  '
$ws.Range('D15').Value = '''92.82'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').Value = '30.008.45'
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').Value = '''14.14'
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('D18').Value = '''5.898'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('D19').Value = '''244.91'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').Value = '''0.000007778'
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('D21').Value = '2.155.53'
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '''0.9997'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').Value = '''0.9996'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = '''6.963'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('D25').Value = '''0.1618'
$ws.Range('E25').Value = '  +25.41%  '
$ws.Range('D26').Value = '''169.27'
$ws.Range('E26').Value = '  +1.95%  '
$ws.Range('D27').Value = '''9.239'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').Value = '''18.94'
$ws.Range('E28').Value = '  +1.22%  '
$ws.Range('D29').Value = '''2.088'
$ws.Range('E29').Value = '  +2.00%  '
$ws.Range('D30').Value = '''1.368'
$ws.Range('E30').Value = '  -3.33%  '
$ws.Range('D31').Value = '''1.516'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('D33').Value = '''0.05610'
$ws.Range('E33').Value = '  +7.27%  '
$ws.Range('D34').Value = '''4.089'
$ws.Range('E34').Value = '  +0.39%  '
$ws.Range('D35').Value = '''1.272'
$ws.Range('E35').Value = '  +1.37%  '
$ws.Range('D36').Value = '''0.7338'
$ws.Range('E36').Value = '  +0.84%  '
$ws.Range('D37').Value = '''2.710'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').Value = '''2.791'
$ws.Range('E39').Value = '  +0.35%  '
$ws.Range('D40').Value = '''0.4433'
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('D41').Value = '''72.46'
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('D42').Value = '''5.982'
$ws.Range('E42').Value = '  -2.93%  '
$ws.Range('D43').Value = '''0.8419'
$ws.Range('E43').Value = '  +1.85%  '
$ws.Range('D44').Value = '''0.9995'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').Value = '''1.894'
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '''7.600'
$ws.Range('E46').Value = '  -0.57%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '''101.04'
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('D48').Value = '''9.733'
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range('D49').Value = '''987.37'
$ws.Range('E49').Value = '  +9.66%  '
$ws.Range('D50').Value = '2.062.25'
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('D51').Value = '''36.25'
